{"js": "// Office.js (Word JavaScript API) script.\n//\n// The questionnaire body has a fixed sequence of 30 paragraphs; 7 of them\n// (the 'phishing message' bodies under the Single/Paired Message sections)\n// get their wording swapped out for new sample messages. Paragraph indices\n// below are 0-based positions in context.document.body.paragraphs, which are\n// stable because no paragraphs are added/removed by this edit.\n//\n// Each replacement is applied with paragraph.insertOoxml(..., replace) instead\n// of clear()/insertText()/insertBreak() so the emitted <w:r>/<w:t>/<w:br/>\n// structure -- including exactly which <w:t> elements need\n// xml:space=\"preserve\" -- matches the target OOXML verbatim.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\n// index -> OOXML for the replacement <w:p>, wrapped in the minimal\n// pkg:package envelope insertOoxml expects.\nconst replacements = [\n  // [7] Single Q1 'Phishing message:' -- was \"Dear Gary  Leal ...\" (fitness tracker)\n  [7, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Want to watch the UEFA league final in HD from the comfort of your own home?</w:t><w:br/><w:br/><w:t>Click the link below for more information</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [10] Single Q2 'Phishing message:' -- was \"Hi Stevie, ...\" (gym induction)\n  [10, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t xml:space=\\\"preserve\\\">Hi Kim Young, your child needs to complete payment and registration for the upcoming school trip. Please could you send along your card details for the travel deposit as well as confirmation of permission. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [13] Single Q3 'Phishing message:' -- was \"Subject: Beauty Offer ... Evelyn\" (beauty offer)\n  [13, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t xml:space=\\\"preserve\\\">Hi Mr Ibarra, </w:t><w:br/><w:br/><w:t>I found you on the DIY website that you signed up to!  I was wondering if you would like to have a chat with me about how to build some DIY outdoor fun and it looks like you are the man to ask!  Click on the link below to chat!</w:t><w:br/><w:br/><w:t>I look forward to hearing from you!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [17] Paired Q1 'Phishing message A:' -- was \"Hello Claude, ...\" (gym payment)\n  [17, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Dear Nicole,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. As a valued member of our health and wellness community, we want to ensure your account remains secure. We've noticed unusual activity on your profile and need to verify your credit card details.</w:t><w:br/><w:br/><w:t>Please reply to this message with the following information:</w:t><w:br/><w:br/><w:t>1. Full Name:</w:t><w:br/><w:t>2. Credit Card Number:</w:t><w:br/><w:t>3. Expiration Date:</w:t><w:br/><w:t>4. CVV Code:</w:t><w:br/><w:br/><w:t>Your security is our top priority. We appreciate your prompt response to help us safeguard your account.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>The Health &amp; Wellness Team</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [19] Paired Q1 'Phishing message B:' -- was \"Subject: Special Offer from Amazon ...\"\n  [19, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Dear Ms. Morrow,</w:t><w:br/><w:br/><w:t>Your package is being held at customs. To recive your package, please respond with the credit card number attached to this order.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [24] Paired Q2 'Phishing message A:' -- was \"Dear customer ...\" (gym membership win)\n  [24, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Subject: Justin, Your Subscription is Expiring Soon</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Dear Justin,</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> We hope this message finds you well. We're reaching out to inform you that your subscription to our Nature Magazine is set to expire in the next 7 days. </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> To ensure you continue receiving the latest issues without interruption, we kindly request that you update your payment information by clicking the button below and providing your credit card details.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> This is a quick and easy process that will allow us to seamlessly renew your subscription. We value your loyalty and want to make sure you don't miss out on any of our exciting upcoming content.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Please click here to update your payment details:</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> [Button: Update Payment Information]</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Thank you for your continued support. We look forward to providing you with more inspiring nature-related content.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Best regards,</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Nature Magazine Subscription Team</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n  // [26] Paired Q2 'Phishing message B:' -- was \"   Dear Esther Lindsley, ...\" (bank of Romania)\n  [26, \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Dear Mr. Miller</w:t><w:br/><w:br/><w:t>Your account is on hold, we're having some trouble with your current billing information. We'll try again, but in th meantime you want to update your payment information</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"],\n];\n\nfor (const [idx, ooxml] of replacements) {\n  paragraphs.items[idx].insertOoxml(ooxml, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Replaces the body text of 7 'phishing message' paragraphs in the\n# questionnaire with new message content, per the commit diff.\n#\n# For each target paragraph we clear the paragraph's Range (excluding the\n# trailing paragraph mark, so paragraph count/formatting stay intact) and\n# then use Range.InsertXML with a literal <w:p> fragment so the resulting\n# run/break structure -- including xml:space=\"preserve\" placement -- matches\n# the target OOXML exactly (letting Word infer xml:space itself can drop the\n# attribute on interior runs that start/end with whitespace).\n\n$d = $word.ActiveDocument\n\n$pkgPrefix = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfunction Set-ParagraphXml($paraIndex, $pFragment) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $rng = $p.Range\n    # Exclude the trailing paragraph mark so the paragraph itself is kept.\n    $rng.End = $rng.End - 1\n    $rng.Text = ''\n    $full = $pkgPrefix + $pFragment + $pkgSuffix\n    [void]$rng.InsertXML($full)\n}\n\n$p7 = '<w:p><w:r><w:t>Want to watch the UEFA league final in HD from the comfort of your own home?</w:t><w:br/><w:br/><w:t>Click the link below for more information</w:t></w:r></w:p>'\nSet-ParagraphXml 8 $p7\n\n$p10 = '<w:p><w:r><w:t xml:space=\"preserve\">Hi Kim Young, your child needs to complete payment and registration for the upcoming school trip. Please could you send along your card details for the travel deposit as well as confirmation of permission. </w:t></w:r></w:p>'\nSet-ParagraphXml 11 $p10\n\n$p13 = '<w:p><w:r><w:t xml:space=\"preserve\">Hi Mr Ibarra, </w:t><w:br/><w:br/><w:t>I found you on the DIY website that you signed up to!  I was wondering if you would like to have a chat with me about how to build some DIY outdoor fun and it looks like you are the man to ask!  Click on the link below to chat!</w:t><w:br/><w:br/><w:t>I look forward to hearing from you!</w:t></w:r></w:p>'\nSet-ParagraphXml 14 $p13\n\n$p17 = '<w:p><w:r><w:t>Dear Nicole,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. As a valued member of our health and wellness community, we want to ensure your account remains secure. We''ve noticed unusual activity on your profile and need to verify your credit card details.</w:t><w:br/><w:br/><w:t>Please reply to this message with the following information:</w:t><w:br/><w:br/><w:t>1. Full Name:</w:t><w:br/><w:t>2. Credit Card Number:</w:t><w:br/><w:t>3. Expiration Date:</w:t><w:br/><w:t>4. CVV Code:</w:t><w:br/><w:br/><w:t>Your security is our top priority. We appreciate your prompt response to help us safeguard your account.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>The Health &amp; Wellness Team</w:t></w:r></w:p>'\nSet-ParagraphXml 18 $p17\n\n$p19 = '<w:p><w:r><w:t>Dear Ms. Morrow,</w:t><w:br/><w:br/><w:t>Your package is being held at customs. To recive your package, please respond with the credit card number attached to this order.</w:t></w:r></w:p>'\nSet-ParagraphXml 20 $p19\n\n$p24 = '<w:p><w:r><w:t>Subject: Justin, Your Subscription is Expiring Soon</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Dear Justin,</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We hope this message finds you well. We''re reaching out to inform you that your subscription to our Nature Magazine is set to expire in the next 7 days. </w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> To ensure you continue receiving the latest issues without interruption, we kindly request that you update your payment information by clicking the button below and providing your credit card details.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> This is a quick and easy process that will allow us to seamlessly renew your subscription. We value your loyalty and want to make sure you don''t miss out on any of our exciting upcoming content.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Please click here to update your payment details:</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> [Button: Update Payment Information]</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Thank you for your continued support. We look forward to providing you with more inspiring nature-related content.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Best regards,</w:t><w:br/><w:t xml:space=\"preserve\"> Nature Magazine Subscription Team</w:t></w:r></w:p>'\nSet-ParagraphXml 25 $p24\n\n$p26 = '<w:p><w:r><w:t>Dear Mr. Miller</w:t><w:br/><w:br/><w:t>Your account is on hold, we''re having some trouble with your current billing information. We''ll try again, but in th meantime you want to update your payment information</w:t><w:br/></w:r></w:p>'\nSet-ParagraphXml 27 $p26\n\nWrite-Output 'done'"}
